$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: now holds data previously in row 9
$ws.Cells.Item(2, 1).Value = 111470636  # A2
$ws.Cells.Item(2, 2).Value = 94134  # B2
$ws.Cells.Item(2, 4).Value = "NT"  # D2
$ws.Cells.Item(2, 5).Value = 53  # E2
$ws.Cells.Item(2, 6).Value = "Vedtrappmossa"  # F2
$ws.Cells.Item(2, 7).Value = "Crossocalyx hellerianus"  # G2
$ws.Cells.Item(2, 8).Value = "(Nees ex Lindenb.) Meyl."  # H2
$ws.Cells.Item(2, 17).Value = 554457.9939421143  # Q2
$ws.Cells.Item(2, 18).Value = 7003163.892755959  # R2
$ws.Cells.Item(2, 26).Value = "14:41"  # Z2
$ws.Cells.Item(2, 28).Value = "14:41"  # AB2

# Row 3: now holds data previously in row 13
$ws.Cells.Item(3, 1).Value = 111470743  # A3
$ws.Cells.Item(3, 2).Value = 78611  # B3
$ws.Cells.Item(3, 4).Value = "LC"  # D3
$ws.Cells.Item(3, 5).Value = 6463  # E3
$ws.Cells.Item(3, 6).Value = "Bårdlav"  # F3
$ws.Cells.Item(3, 7).Value = "Nephroma parile"  # G3
$ws.Cells.Item(3, 17).Value = 554457.9939421143  # Q3
$ws.Cells.Item(3, 18).Value = 7003163.892755959  # R3
$ws.Cells.Item(3, 26).Value = "14:41"  # Z3
$ws.Cells.Item(3, 28).Value = "14:41"  # AB3
$ws.Cells.Item(3, 29).ClearContents()  # AC3

# Row 4: now holds data previously in row 6
$ws.Cells.Item(4, 1).Value = 111470448  # A4
$ws.Cells.Item(4, 17).Value = 554488.5866359913  # Q4
$ws.Cells.Item(4, 18).Value = 7003175.257923778  # R4
$ws.Cells.Item(4, 19).Value = 22  # S4
$ws.Cells.Item(4, 26).Value = "14:59"  # Z4
$ws.Cells.Item(4, 28).Value = "14:59"  # AB4

# Row 5: now holds data previously in row 10
$ws.Cells.Item(5, 1).Value = 111471685  # A5
$ws.Cells.Item(5, 17).Value = 554595.0694405095  # Q5
$ws.Cells.Item(5, 18).Value = 7003142.694495555  # R5
$ws.Cells.Item(5, 26).Value = "15:49"  # Z5
$ws.Cells.Item(5, 28).Value = "15:49"  # AB5
$ws.Cells.Item(5, 29).ClearContents()  # AC5

# Row 6: now holds data previously in row 7
$ws.Cells.Item(6, 1).Value = 111470486  # A6
$ws.Cells.Item(6, 2).Value = 78578  # B6
$ws.Cells.Item(6, 4).Value = "NT"  # D6
$ws.Cells.Item(6, 5).Value = 6458  # E6
$ws.Cells.Item(6, 6).Value = "Lunglav"  # F6
$ws.Cells.Item(6, 7).Value = "Lobaria pulmonaria"  # G6
$ws.Cells.Item(6, 8).Value = "(L.) Hoffm."  # H6
$ws.Cells.Item(6, 26).Value = "14:41"  # Z6
$ws.Cells.Item(6, 28).Value = "14:41"  # AB6

# Row 7: now holds data previously in row 3
$ws.Cells.Item(7, 1).Value = 111471797  # A7
$ws.Cells.Item(7, 2).Value = 77515  # B7
$ws.Cells.Item(7, 5).Value = 6425  # E7
$ws.Cells.Item(7, 6).Value = "Garnlav"  # F7
$ws.Cells.Item(7, 7).Value = "Alectoria sarmentosa"  # G7
$ws.Cells.Item(7, 8).Value = "(Ach.) Ach."  # H7
$ws.Cells.Item(7, 17).Value = 554597.2688619854  # Q7
$ws.Cells.Item(7, 18).Value = 7003280.616068945  # R7
$ws.Cells.Item(7, 19).Value = 25  # S7
$ws.Cells.Item(7, 26).Value = "15:49"  # Z7
$ws.Cells.Item(7, 28).Value = "15:49"  # AB7
$ws.Cells.Item(7, 29).Value = "På tall"  # AC7

# Row 8: now holds data previously in row 2
$ws.Cells.Item(8, 1).Value = 111471083  # A8
$ws.Cells.Item(8, 17).Value = 554499.1143642976  # Q8
$ws.Cells.Item(8, 18).Value = 7003141.52872613  # R8
$ws.Cells.Item(8, 26).Value = "15:31"  # Z8
$ws.Cells.Item(8, 28).Value = "15:31"  # AB8

# Row 9: now holds data previously in row 5
$ws.Cells.Item(9, 1).Value = 111469986  # A9
$ws.Cells.Item(9, 2).Value = 77515  # B9
$ws.Cells.Item(9, 5).Value = 6425  # E9
$ws.Cells.Item(9, 6).Value = "Garnlav"  # F9
$ws.Cells.Item(9, 7).Value = "Alectoria sarmentosa"  # G9
$ws.Cells.Item(9, 8).Value = "(Ach.) Ach."  # H9
$ws.Cells.Item(9, 17).Value = 554489.6113782075  # Q9
$ws.Cells.Item(9, 18).Value = 7003329.432399829  # R9
$ws.Cells.Item(9, 26).Value = "00:00"  # Z9
$ws.Cells.Item(9, 28).Value = "00:00"  # AB9
$ws.Cells.Item(9, 29).Value = "Rikligt på tall"  # AC9

# Row 10: now holds data previously in row 4
$ws.Cells.Item(10, 1).Value = 111470792  # A10
$ws.Cells.Item(10, 2).Value = 96348  # B10
$ws.Cells.Item(10, 4).Value = "VU"  # D10
$ws.Cells.Item(10, 5).Value = 220787  # E10
$ws.Cells.Item(10, 6).Value = "Knärot"  # F10
$ws.Cells.Item(10, 7).Value = "Goodyera repens"  # G10
$ws.Cells.Item(10, 8).Value = "(L.) R. Br."  # H10
$ws.Cells.Item(10, 17).Value = 554440.9784625648  # Q10
$ws.Cells.Item(10, 18).Value = 7003152.756292564  # R10
$ws.Cells.Item(10, 26).Value = "15:19"  # Z10
$ws.Cells.Item(10, 28).Value = "15:19"  # AB10

# Row 12: now holds data previously in row 8
$ws.Cells.Item(12, 1).Value = 111470245  # A12
$ws.Cells.Item(12, 2).Value = 96348  # B12
$ws.Cells.Item(12, 4).Value = "VU"  # D12
$ws.Cells.Item(12, 5).Value = 220787  # E12
$ws.Cells.Item(12, 6).Value = "Knärot"  # F12
$ws.Cells.Item(12, 7).Value = "Goodyera repens"  # G12
$ws.Cells.Item(12, 8).Value = "(L.) R. Br."  # H12
$ws.Cells.Item(12, 17).Value = 554481.1995954363  # Q12
$ws.Cells.Item(12, 18).Value = 7003291.317192273  # R12

# Row 13: now holds data previously in row 12
$ws.Cells.Item(13, 1).Value = 111470685  # A13
$ws.Cells.Item(13, 2).Value = 77267  # B13
$ws.Cells.Item(13, 4).Value = "NT"  # D13
$ws.Cells.Item(13, 5).Value = 6446  # E13
$ws.Cells.Item(13, 6).Value = "Kolflarnlav"  # F13
$ws.Cells.Item(13, 7).Value = "Carbonicola anthracophila"  # G13
$ws.Cells.Item(13, 8).Value = "(Nyl.) Bendiksby & Timdal"  # H13
